$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at M. This pushes the old (empty) M -> N, the old
#    "Time" helper block N,O,P (labels/values/formulas) -> O,P,Q, and Excel
#    auto-rewrites every formula that referenced the old O/O3 cells so they
#    now point at the new P/P3 cells.
$ws.Columns("M").Insert()

# 2. Header row: new "Top Speed"/"Min Speed" columns.
$ws.Range("L1").Value2 = "Top Speed (px/t)"
$ws.Range("M1").Value2 = "Min Speed (px/t)"

# Give the (still blank) M1/N1 header cells the same centered style used by
# the rest of row 1 / the time block.
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("N1").HorizontalAlignment = -4108

# 3. Re-point the "Days" value cell's alignment: after the shift it is s="2"
#    (font/border/center) - the target only wants the plain center style.
$ws.Range("P2").HorizontalAlignment = -4108

# 4. Top speed becomes relative to the Hawk's (row 6) top speed.
$ws.Range("L2").Formula = "=L6*0.68"
$ws.Range("L3").Formula = "=L6*0.64"
$ws.Range("L4").Formula = "=L6*0.16"
$ws.Range("L5").Formula = "=L6*0.29"
$ws.Range("L6").Value2 = 10
$ws.Range("L7").Formula = "=L6*0.4"
$ws.Range("L8").Formula = "=L6*0.4"
$ws.Range("L9").Formula = "=L6*0.15"
$ws.Range("L10").Formula = "=L6*0.128"

# 5. New "Min Speed" column = half the top speed.
$ws.Range("M2").Formula = "=L2/2"
$ws.Range("M3:M10").Formula = "=L3/2"
# These cells keep the workbook's default (unstyled) formatting, unlike the
# rest of the table.
$ws.Range("M2:M10").ClearFormats()

# 6. Re-create the shared formula group for the "seconds per tick" column
#    (previously P3:P5, now Q3:Q5) so it serialises the same way it used to.
$ws.Range("Q3:Q5").Formula = "=P3/30"

# 7. Column widths for the two newly introduced columns.
$ws.Range("M1:N1").ColumnWidth = 16.59

# 8. Alternate-row banding via conditional formatting.
$cfRange = $ws.Range("A1:XFD10")
$cf = $cfRange.FormatConditions.Add(2, 0, '"MOD(ROW(),2)=1"')
$cf.Interior.Pattern = 1
$cf.Interior.Color = 5287936

# 9. View state: reset zoom and move the selection.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("I9").Select()
